$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "194088"
$ws.Range("E3").NumberFormat = "General"
